# Fill column B ("Values") on each color-named worksheet with the
# worksheet's own name, mirroring the per-sheet "color" value.
# On the "D Green" sheet only, rows 7, 8 and 12 (BMW X2, Renault Rafale HEV,
# Renault Espace) are left blank / untouched.

$wb = $excel.ActiveWorkbook

$skipRowsBySheet = @{
    "D Green" = @(7, 8, 12)
}

foreach ($ws in $wb.Worksheets) {
    $sheetName = $ws.Name

    $skipRows = @()
    if ($skipRowsBySheet.ContainsKey($sheetName)) {
        $skipRows = $skipRowsBySheet[$sheetName]
    }

    for ($row = 2; $row -le 15; $row++) {
        if ($skipRows -contains $row) {
            $ws.Cells.Item($row, 2).ClearContents()
            continue
        }
        $ws.Cells.Item($row, 2).Value = $sheetName
    }
}
